$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cria o cabecalho na primeira linha (A1:B1)
$headers = @("a", "b")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# Dados existentes deslocados/preenchidos na linha 2 (A2:B2)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
